# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 3.123541145015474
$ws.Range("E2").Value = 2.651224325850143
$ws.Range("C3").Value = 3.959010658874851
$ws.Range("E3").Value = 2.48038747547068
$ws.Range("C4").Value = 4.722695063536686
$ws.Range("E4").Value = 3.580489688170352
$ws.Range("C5").Value = 6.739021039846627
$ws.Range("E5").Value = 2.457570659067509
$ws.Range("C6").Value = 2.619839412265601
$ws.Range("E6").Value = 1.929582042845213
$ws.Range("C7").Value = -0.7919564768266385
$ws.Range("E7").Value = 1.978648203842193
$ws.Range("C8").Value = 1.877689851450803
$ws.Range("E8").Value = 3.188468414048606
$ws.Range("C9").Value = 2.705004599189187
$ws.Range("E9").Value = 2.471779821159181
$ws.Range("C10").Value = 1.110374544249249
$ws.Range("E10").Value = 2.958768964947134
$ws.Range("C11").Value = 2.267566233338814
$ws.Range("E11").Value = 2.824054578526636
$ws.Range("C12").Value = 2.688433258834588
$ws.Range("E12").Value = 2.558605038029849
$ws.Range("C13").Value = 1.014079695989589
$ws.Range("E13").Value = 2.1961756525541
$ws.Range("C14").Value = 3.013853578092252
$ws.Range("E14").Value = 2.721358666668894
$ws.Range("C15").Value = 1.331295149770684
$ws.Range("E15").Value = 1.104792173470215
$ws.Range("C16").Value = 0.04589006555719699
$ws.Range("E16").Value = 1.367746832546346
$ws.Range("C17").Value = 0.009546395482029624
$ws.Range("E17").Value = 0.9083591402091473
$ws.Range("C18").Value = 0.8709390141433015
$ws.Range("E18").Value = 1.177606443599988
$ws.Range("C19").Value = 0.7652063367885598
$ws.Range("E19").Value = 1.680477107215861
$ws.Range("C20").Value = 2.267579219134386
$ws.Range("E20").Value = 2.225111592343887
$ws.Range("C21").Value = 3.146753122914103
$ws.Range("E21").Value = 1.799793194919874
$ws.Range("C22").Value = 1.769033835366818
$ws.Range("E22").Value = 0.2856860139923256
$ws.Range("C23").Value = -4.774715709990263
$ws.Range("E23").Value = 0.754926127539246
$ws.Range("C24").Value = 1.95493704440024
$ws.Range("E24").Value = 2.928378677701393
$ws.Range("C25").Value = 3.478075069442799
$ws.Range("E25").Value = 1.768431385360159
$ws.Range("C26").Value = 1.232342134690434
$ws.Range("E26").Value = 1.134779475590464
$ws.Range("C27").Value = 0.2542811494408159
$ws.Range("E27").Value = 1.353526127153426
$ws.Range("C28").Value = 1.519778766382096
$ws.Range("E28").Value = 0.7212678493511149
$ws.Range("C29").Value = 1.469441753880329
$ws.Range("E29").Value = 1.350400980499855
$ws.Range("C30").Value = 1.638203081492495
$ws.Range("E30").Value = 1.323745783269614
$ws.Range("C31").Value = 2.268697431234346
$ws.Range("E31").Value = 2.455413743911294
$ws.Range("C32").Value = 1.984425467899631
$ws.Range("E32").Value = 0.6687400825358569
$ws.Range("C33").Value = 0.6066448776129052
$ws.Range("E33").Value = 0.6230021429014077
$ws.Range("C34").Value = -4.243076347305386
$ws.Range("E34").Value = -2.163103471150829
$ws.Range("C35").Value = 1.438499295329754
$ws.Range("E35").Value = 1.882119284761474
$ws.Range("C36").Value = 1.906593537051537
$ws.Range("E36").Value = 1.473394465200051
$ws.Range("C37").Value = 0.08348019664223827
$ws.Range("E37").Value = 0.8702074629614476
$ws.Range("C38").Value = -0.214505326882275
$ws.Range("E38").Value = 0.8092352694139215
$ws.Range("C39").Value = 0.1651547428133782
$ws.Range("E39").Value = 0.7536567386490001
